$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BY is column index 77. The sheet's "group" columns are extended
# with 13 more columns, BZ..CL (indexes 78..90). Every new cell in a given
# row repeats whatever value (and formatting) its row already carries in
# column BY - including row 10, whose BY cell is blank.
for ($r = 2; $r -le 15; $r++) {
    $srcCell = $ws.Cells.Item($r, 77)
    $srcVal = $srcCell.Text
    for ($c = 78; $c -le 90; $c++) {
        $dstCell = $ws.Cells.Item($r, $c)
        $dstCell.Value = $srcVal
        $dstCell.Style = $srcCell.Style
    }
}
